$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- E4: append the new conformance experiment to the "Other experiments" list ---
$ws.Cells.Item(4, 5).Value = "3.1S_100yrControl, 2.1_FutureTimeSlice"

# --- E47:E55: append "; 2.1_FutureTimeSlice" to the shared "other experiments" text ---
$suffix = "; 2.1_FutureTimeSlice"
for ($r = 47; $r -le 55; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value = $cell.Text + $suffix
}

# --- Row height tweaks on rows 55 and 56 ---
$ws.Rows.Item(55).RowHeight = 67.5
$ws.Rows.Item(56).RowHeight = 69.75

# --- New rows 84 and 85 describing bc.062 / bc.063 (SST / Sea Ice anomalies) ---
$ws.Cells.Item(84, 2).Value = "bc.062"
$ws.Cells.Item(85, 2).Value = "bc.063"

$ws.Cells.Item(84, 3).Value = "Prescribed SST anomolies (relative to exp 3.3) based on either (a) climatology difference between years 2026-2035 (from expt 4.1) and years 1979-2008) (from expt 3.2), OR (b) climatology difference between years 2026-2035 (from expt 1.2) and years 1979-2008 (from expt1.1)"
$ws.Cells.Item(84, 4).Value = "2.1_FutureTimeSlice"

$ws.Cells.Item(85, 3).Value = "Prescribed Sea Ice anomolies (relative to exp 3.3) based on either (a) climatology difference between years 2026-2035 (from expt 4.1) and years 1979-2008) (from expt 3.2), OR (b) climatology difference between years 2026-2035 (from expt 1.2) and years 1979-2008 (from expt1.1)"
$ws.Cells.Item(85, 4).Value = "2.1_FutureTimeSlice"

# --- Match the styling/format used by the row above (B/D id style, C/E wrapped text style) ---
$ws.Range("B83:E83").Copy()
$ws.Range("B84:E84").PasteSpecial(-4122)
$ws.Range("B85:E85").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(84).RowHeight = 78.75
$ws.Rows.Item(85).RowHeight = 60

# --- Update the view state: scrolled position & active selection ---
$ws.Application.ActiveWindow.ScrollRow = 70
$ws.Range("D87").Select()
